$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update transition-matrix probabilities (recomputed after adding more simulated games)
# Row 2
$ws.Range("B2").Value = 0.2103960396039604
$ws.Range("C2").Value = 0.5297029702970297
$ws.Range("J2").Value = 0.007425742574257425
$ws.Range("P2").Value = 0.150990099009901
$ws.Range("S2").Value = 0.1014851485148515

# Row 3
$ws.Range("B3").Value = 0.00881057268722467
$ws.Range("C3").Value = 0.04405286343612335
$ws.Range("J3").Value = 0.02643171806167401
$ws.Range("P3").Value = 0.7268722466960352
$ws.Range("S3").Value = 0.1938325991189427

# Row 4
$ws.Range("J4").Value = 0.04081632653061224
$ws.Range("P4").Value = 0.6938775510204082
$ws.Range("S4").Value = 0.2653061224489796

# Row 5
$ws.Range("J5").Value = 1

# Row 6
$ws.Range("B6").Value = 0.08085106382978724
$ws.Range("D6").Value = 0.01276595744680851
$ws.Range("F6").Value = 0.1148936170212766
$ws.Range("J6").Value = 0.1872340425531915
$ws.Range("O6").Value = 0.02978723404255319
$ws.Range("Q6").Value = 0.1702127659574468
$ws.Range("R6").Value = 0.07659574468085106
$ws.Range("S6").Value = 0.3276595744680851

# Row 7
$ws.Range("B7").Value = 0.1348837209302326
$ws.Range("D7").Value = 0.009302325581395349
$ws.Range("F7").Value = 0.03720930232558139
$ws.Range("J7").Value = 0.1395348837209302
$ws.Range("O7").Value = 0.01395348837209302
$ws.Range("Q7").Value = 0.1534883720930233
$ws.Range("R7").Value = 0.08837209302325581
$ws.Range("S7").Value = 0.4232558139534884

# Row 8
$ws.Range("B8").Value = 0.1354581673306773
$ws.Range("D8").Value = 0.0298804780876494
$ws.Range("E8").Value = 0.00199203187250996
$ws.Range("F8").Value = 0.06374501992031872
$ws.Range("J8").Value = 0.1294820717131474
$ws.Range("O8").Value = 0.01394422310756972
$ws.Range("Q8").Value = 0.1414342629482072
$ws.Range("R8").Value = 0.09760956175298804
$ws.Range("S8").Value = 0.3864541832669323

# Row 9
$ws.Range("B9").Value = 0.1003861003861004
$ws.Range("D9").Value = 0.0193050193050193
$ws.Range("F9").Value = 0.05019305019305019
$ws.Range("J9").Value = 0.1428571428571428
$ws.Range("O9").Value = 0.03088803088803089
$ws.Range("Q9").Value = 0.1737451737451738
$ws.Range("R9").Value = 0.09266409266409266
$ws.Range("S9").Value = 0.38996138996139

# Row 10
$ws.Range("B10").Value = 0.1239130434782609
$ws.Range("D10").Value = 0.01811594202898551
$ws.Range("E10").Value = 0.0007246376811594203
$ws.Range("F10").Value = 0.06521739130434782
$ws.Range("J10").Value = 0.1166666666666667
$ws.Range("O10").Value = 0.02246376811594203
$ws.Range("Q10").Value = 0.2108695652173913
$ws.Range("R10").Value = 0.0782608695652174
$ws.Range("S10").Value = 0.363768115942029

# Row 11
$ws.Range("G11").Value = 0.1466666666666667
$ws.Range("J11").Value = 0.1093333333333333
$ws.Range("K11").Value = 0.2186666666666667
$ws.Range("L11").Value = 0.512
$ws.Range("S11").Value = 0.01333333333333333

# Row 12
$ws.Range("G12").Value = 0.7193877551020408
$ws.Range("J12").Value = 0.2397959183673469
$ws.Range("K12").Value = 0.01020408163265306
$ws.Range("L12").Value = 0.01530612244897959
$ws.Range("S12").Value = 0.01530612244897959

# Row 13
$ws.Range("G13").Value = 0.6122448979591837
$ws.Range("J13").Value = 0.3265306122448979
$ws.Range("S13").Value = 0.06122448979591837

# Row 15
$ws.Range("F15").Value = 0.016
$ws.Range("H15").Value = 0.148
$ws.Range("I15").Value = 0.08400000000000001
$ws.Range("J15").Value = 0.272
$ws.Range("K15").Value = 0.076
$ws.Range("M15").Value = 0.024
$ws.Range("O15").Value = 0.08799999999999999
$ws.Range("S15").Value = 0.292

# Row 16
$ws.Range("F16").Value = 0.01587301587301587
$ws.Range("H16").Value = 0.2103174603174603
$ws.Range("I16").Value = 0.1150793650793651
$ws.Range("J16").Value = 0.3015873015873016
$ws.Range("K16").Value = 0.1150793650793651
$ws.Range("M16").Value = 0.02380952380952381
$ws.Range("N16").Value = 0.003968253968253968
$ws.Range("O16").Value = 0.04761904761904762
$ws.Range("S16").Value = 0.1666666666666667

# Row 17
$ws.Range("F17").Value = 0.01680672268907563
$ws.Range("H17").Value = 0.157563025210084
$ws.Range("I17").Value = 0.1365546218487395
$ws.Range("J17").Value = 0.3907563025210084
$ws.Range("K17").Value = 0.09663865546218488
$ws.Range("M17").Value = 0.01680672268907563
$ws.Range("O17").Value = 0.06512605042016807
$ws.Range("S17").Value = 0.1197478991596639

# Row 18
$ws.Range("F18").Value = 0.01388888888888889
$ws.Range("H18").Value = 0.1296296296296296
$ws.Range("I18").Value = 0.125
$ws.Range("J18").Value = 0.4537037037037037
$ws.Range("K18").Value = 0.1018518518518518
$ws.Range("M18").Value = 0.02314814814814815
$ws.Range("O18").Value = 0.05555555555555555
$ws.Range("S18").Value = 0.09722222222222222

# Row 19
$ws.Range("F19").Value = 0.01316701316701317
$ws.Range("H19").Value = 0.2141372141372141
$ws.Range("I19").Value = 0.08177408177408177
$ws.Range("J19").Value = 0.358974358974359
$ws.Range("K19").Value = 0.1205821205821206
$ws.Range("M19").Value = 0.01732501732501732
$ws.Range("O19").Value = 0.06098406098406099
$ws.Range("S19").Value = 0.1330561330561331
